$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.135.35'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '2.760.09'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.76'
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.18'
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.605'
$ws.Range("E8").Value = '  -3.27%  '

$ws.Range("E9").Value = '  -1.71%  '

$ws.Range("E10").Value = '  +4.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.388'
$ws.Range("E11").Value = '  -1.57%  '

$ws.Range("E12").Value = '  -16.10%  '

$ws.Range("D13").Value = '3.248.78'
$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.94'
$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").Value = '63.749.06'
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("E16").Value = '  -2.32%  '

$ws.Range("D17").Value = '2.764.30'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.14'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("E19").Value = '  -1.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '357.58'
$ws.Range("E20").Value = '  -1.79%  '

$ws.Range("E21").Value = '  -3.99%  '

$ws.Range("E22").Value = '  +0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.535'
$ws.Range("E23").Value = '  -1.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.55'
$ws.Range("E24").Value = '  -2.06%  '

$ws.Range("E25").Value = '  -1.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.61'
$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '0.0₃0911'
$ws.Range("E28").Value = '  -1.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.32'
$ws.Range("E29").Value = '  +0.88%  '

$ws.Range("E30").Value = '  -3.21%  '

$ws.Range("E31").Value = '  -0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.81'
$ws.Range("E32").Value = '  -2.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.30'
$ws.Range("E33").Value = '  -1.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.94'
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").Value = '  +2.00%  '

$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("E37").Value = '  -0.50%  '

$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.37'
$ws.Range("E39").Value = '  +2.41%  '

$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '345.20'
$ws.Range("E40").Value = '  +1.92%  '

$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.23'
$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.46'
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.87'
$ws.Range("E44").Value = '  -2.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0592'
$ws.Range("E45").Value = '  -2.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0256'
$ws.Range("E46").Value = '  -1.35%  '

$ws.Range("E47").Value = '  -2.03%  '

$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("E49").Value = '  -1.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.06'
$ws.Range("E51").Value = '  +0.11%  '
